$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.399.04"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "1.843.38"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("D4").Value = "'0.9986"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'240.22"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "'0.6323"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("D7").Value = "'0.9997"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2905"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'25.02"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("D11").Value = "'0.07743"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").Value = "1.846.12"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").Value = "'4.988"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "'0.6788"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "'82.02"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "'6.270"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +2.90%  "

$ws.Range("D18").Value = "29.462.88"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").Value = "'229.92"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "'12.33"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").Value = "'0.9995"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").Value = "'7.420"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D24").Value = "'158.33"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("D25").Value = "'8.498"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("D26").Value = "'0.1356"
$ws.Range("D26").NumberFormat = "General"

$ws.Range("D27").Value = "'17.46"
$ws.Range("D27").NumberFormat = "General"

$ws.Range("D28").Value = "'0.06551"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +15.43%  "

$ws.Range("E29").Value = "  +2.21%  "

$ws.Range("D30").Value = "'1.485"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("D31").Value = "'4.077"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("D32").Value = "'4.053"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").Value = "'1.839"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").Value = "'1.141"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").Value = "'0.6986"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").Value = "'0.01859"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.60%  "

$ws.Range("D38").Value = "1.250.61"
$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").Value = "'6.776"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("D41").Value = "'0.9345"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +3.44%  "

$ws.Range("D42").Value = "'0.9999"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "2.004.85"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'101.00"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "'65.48"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("E46").Value = "  +4.04%  "

$ws.Range("D47").Value = "'7.063"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("E48").Value = "  +3.92%  "

$ws.Range("D49").Value = "'9.017"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").Value = "'0.1148"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").Value = "'0.3910"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.62%  "

